$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1

$wb.Save()
